$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the date value in A2 and the hyperlink-text value in B2, leaving the
# cell formatting/styles intact (empty cells, same style indices).
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# Remove the hyperlink attached to B2 (drops <hyperlinks> and its relationship).
$ws.Hyperlinks.Delete()

# Selection moves from B8 to the whole of row 2 (active cell A2).
$ws.Range("A2:XFD2").Select()
